$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "want to go" counts in column F, rows 2-7
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 13668
$wsExhibit.Range("F3").Value = 323
$wsExhibit.Range("F4").Value = 662
$wsExhibit.Range("F5").Value = 229
$wsExhibit.Range("F6").Value = 479
$wsExhibit.Range("F7").Value = 1394

# Sheet "全部类型" (All types) - same underlying events, different row layout
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 13668
$wsAll.Range("F3").Value = 323
$wsAll.Range("F4").Value = 662
$wsAll.Range("F5").Value = 229
$wsAll.Range("F8").Value = 479
$wsAll.Range("F9").Value = 1394
